$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.84187548871015
$ws.Range("C2").Value = 8.925680154868447
$ws.Range("D2").Value = 6.204353910365369
$ws.Range("F2").Value = 31.14566759223453
$ws.Range("G2").Value = 40.74255314430381
$ws.Range("H2").Value = 17.18685757625622
$ws.Range("K2").Value = 10.79395292186509
$ws.Range("L2").Value = 10.79280397497298
$ws.Range("N2").Value = 20.56206914240566
$ws.Range("B3").Value = 14.60693785875613
$ws.Range("C3").Value = 8.915190186543249
$ws.Range("D3").Value = 6.191737039343295
$ws.Range("F3").Value = 31.11291363922697
$ws.Range("G3").Value = 40.66796731271003
$ws.Range("H3").Value = 17.21944112290621
$ws.Range("K3").Value = 10.63363062840116
$ws.Range("L3").Value = 10.76875351809058
$ws.Range("N3").Value = 20.62376452128861
$ws.Range("B4").Value = 14.46456576195212
$ws.Range("C4").Value = 8.908679823197051
$ws.Range("D4").Value = 6.183792707589295
$ws.Range("F4").Value = 31.10045249013235
$ws.Range("G4").Value = 40.63339783930886
$ws.Range("H4").Value = 17.24234776796768
$ws.Range("K4").Value = 10.53655150018425
$ws.Range("L4").Value = 10.75618725647132
$ws.Range("N4").Value = 20.66351643233888
$ws.Range("B5").Value = 14.40709949539428
$ws.Range("C5").Value = 8.906009644647501
$ws.Range("D5").Value = 6.180506079593889
$ws.Range("F5").Value = 31.09730112725395
$ws.Range("G5").Value = 40.62214026960328
$ws.Range("H5").Value = 17.25241085408757
$ws.Range("K5").Value = 10.49738417605723
$ws.Range("L5").Value = 10.7516233789008
$ws.Range("N5").Value = 20.68018710149367
$ws.Range("B6").Value = 14.39759286378525
$ws.Range("C6").Value = 8.905565236861165
$ws.Range("D6").Value = 6.179957375369755
$ws.Range("F6").Value = 31.09689426957456
$ws.Range("G6").Value = 40.62044201112413
$ws.Range("H6").Value = 17.25412579161756
$ws.Range("K6").Value = 10.49090571245825
$ws.Range("L6").Value = 10.75089928483071
$ws.Range("N6").Value = 20.68298375810575
$ws.Range("B7").Value = 14.46378841519193
$ws.Range("C7").Value = 8.908643881262579
$ws.Range("D7").Value = 6.183748581834224
$ws.Range("F7").Value = 31.10040218598258
$ws.Range("G7").Value = 40.63323455142858
$ws.Range("H7").Value = 17.24248053403745
$ws.Range("K7").Value = 10.53602161634684
$ws.Range("L7").Value = 10.75612344690957
$ws.Range("N7").Value = 20.66373934859674
$ws.Range("B8").Value = 14.76052596927824
$ws.Range("C8").Value = 8.922077823798457
$ws.Range("D8").Value = 6.200044858089681
$ws.Range("F8").Value = 31.13278855299655
$ws.Range("G8").Value = 40.71451080277586
$ws.Range("H8").Value = 17.19748996324197
$ws.Range("K8").Value = 10.73842258672073
$ws.Range("L8").Value = 10.78405688197884
$ws.Range("N8").Value = 20.58295407376732
$ws.Range("B9").Value = 15.35388871194699
$ws.Range("C9").Value = 8.947857315601032
$ws.Range("D9").Value = 6.23041706769108
$ws.Range("F9").Value = 31.25679932803057
$ws.Range("G9").Value = 40.96254563384804
$ws.Range("H9").Value = 17.13231061184141
$ws.Range("K9").Value = 11.14384850105221
$ws.Range("L9").Value = 10.85611907876538
$ws.Range("N9").Value = 20.43932919709253
$ws.Range("B10").Value = 15.79222920875968
$ws.Range("C10").Value = 8.966437067224206
$ws.Range("D10").Value = 6.251743373872209
$ws.Range("F10").Value = 31.38443435710543
$ws.Range("G10").Value = 41.19806741823677
$ws.Range("H10").Value = 17.0985176773932
$ws.Range("K10").Value = 11.44388668913246
$ws.Range("L10").Value = 10.91933218112379
$ws.Range("N10").Value = 20.34275838429203
$ws.Range("B11").Value = 15.99121757978611
$ws.Range("C11").Value = 8.974807768416934
$ws.Range("D11").Value = 6.261226107537859
$ws.Range("F11").Value = 31.45032082435332
$ws.Range("G11").Value = 41.31656304190921
$ws.Range("H11").Value = 17.08621232485908
$ws.Range("K11").Value = 11.58022770311338
$ws.Range("L11").Value = 10.9502521953933
$ws.Range("N11").Value = 20.30075475350751
$ws.Range("B12").Value = 16.06643387856702
$ws.Range("C12").Value = 8.977965514397093
$ws.Range("D12").Value = 6.264785105801959
$ws.Range("F12").Value = 31.47638368128206
$ws.Range("G12").Value = 41.36304396414802
$ws.Range("H12").Value = 17.08199407212559
$ws.Range("K12").Value = 11.63178474255671
$ws.Range("L12").Value = 10.96226556777566
$ws.Range("N12").Value = 20.28512512267212
$ws.Range("B13").Value = 16.05024195298142
$ws.Range("C13").Value = 8.977285979334082
$ws.Range("D13").Value = 6.264020038172897
$ws.Range("F13").Value = 31.47072128298668
$ws.Range("G13").Value = 41.35296230282277
$ws.Range("H13").Value = 17.08288290577871
$ws.Range("K13").Value = 11.62068502421049
$ws.Range("L13").Value = 10.95966483340878
$ws.Range("N13").Value = 20.28847897123588
$ws.Range("B14").Value = 15.99740882977879
$ws.Range("C14").Value = 8.975067801802467
$ws.Range("D14").Value = 6.261519550543484
$ws.Range("F14").Value = 31.45244279819546
$ws.Range("G14").Value = 41.32035495010965
$ws.Range("H14").Value = 17.08585643583505
$ws.Range("K14").Value = 11.58447106816188
$ws.Range("L14").Value = 10.95123447654437
$ws.Range("N14").Value = 20.29946336309975
$ws.Range("B15").Value = 15.96502704389648
$ws.Range("C15").Value = 8.973707523049391
$ws.Range("D15").Value = 6.2599837608952
$ws.Range("F15").Value = 31.44139126940114
$ws.Range("G15").Value = 41.30059082003645
$ws.Range("H15").Value = 17.08773531884295
$ws.Range("K15").Value = 11.56227806613348
$ws.Range("L15").Value = 10.94611010988777
$ws.Range("N15").Value = 20.30622756882115
$ws.Range("B16").Value = 15.77920988807991
$ws.Range("C16").Value = 8.96588835483605
$ws.Range("D16").Value = 6.251119214803215
$ws.Range("F16").Value = 31.38028497496848
$ws.Range("G16").Value = 41.19054998917713
$ws.Range("H16").Value = 17.09938362306964
$ws.Range("K16").Value = 11.43496910331299
$ws.Range("L16").Value = 10.91735448158571
$ws.Range("N16").Value = 20.34554212127617
$ws.Range("B17").Value = 15.66505542057061
$ws.Range("C17").Value = 8.961070462400265
$ws.Range("D17").Value = 6.245624760135125
$ws.Range("F17").Value = 31.34479421069919
$ws.Range("G17").Value = 41.12593594198621
$ws.Range("H17").Value = 17.1073154311699
$ws.Range("K17").Value = 11.35679444416322
$ws.Range("L17").Value = 10.9002635791016
$ws.Range("N17").Value = 20.37015324087101
$ws.Range("B18").Value = 15.59936215387493
$ws.Range("C18").Value = 8.95829170561473
$ws.Range("D18").Value = 6.242443936247795
$ws.Range("F18").Value = 31.32511817798564
$ws.Range("G18").Value = 41.08984208950569
$ws.Range("H18").Value = 17.11216627920531
$ws.Range("K18").Value = 11.31181948965368
$ws.Range("L18").Value = 10.89063745631911
$ws.Range("N18").Value = 20.38449032635722
$ws.Range("B19").Value = 15.57711597389454
$ws.Range("C19").Value = 8.957349566455886
$ws.Range("D19").Value = 6.241363446900599
$ws.Range("F19").Value = 31.31858318937489
$ws.Range("G19").Value = 41.07780585770463
$ws.Range("H19").Value = 17.11385825860028
$ws.Range("K19").Value = 11.29659144419809
$ws.Range("L19").Value = 10.88741346549656
$ws.Range("N19").Value = 20.38937580268579
$ws.Range("B20").Value = 15.67721150718644
$ws.Range("C20").Value = 8.961584127835099
$ws.Range("D20").Value = 6.246211785289067
$ws.Range("F20").Value = 31.34849603634611
$ws.Range("G20").Value = 41.13270359428017
$ws.Range("H20").Value = 17.10644119547385
$ws.Range("K20").Value = 11.36511777954559
$ws.Range("L20").Value = 10.90206185607907
$ws.Range("N20").Value = 20.36751457507101
$ws.Range("B21").Value = 16.01293149361258
$ws.Range("C21").Value = 8.975719664435628
$ws.Range("D21").Value = 6.262254873952831
$ws.Range("F21").Value = 31.4577815275293
$ws.Range("G21").Value = 41.32988904928153
$ws.Range("H21").Value = 17.08497105230607
$ws.Range("K21").Value = 11.59511034689228
$ws.Range("L21").Value = 10.95370246343111
$ws.Range("N21").Value = 20.29622949080149
$ws.Range("B22").Value = 16.23151614153889
$ws.Range("C22").Value = 8.984887608326362
$ws.Range("D22").Value = 6.272553770236948
$ws.Range("F22").Value = 31.53568712096689
$ws.Range("G22").Value = 41.46812789862395
$ws.Range("H22").Value = 17.07351269743915
$ws.Range("K22").Value = 11.74497967223779
$ws.Range("L22").Value = 10.98922488254108
$ws.Range("N22").Value = 20.25125027768821
$ws.Range("B23").Value = 16.11495369298503
$ws.Range("C23").Value = 8.980001070359419
$ws.Range("D23").Value = 6.267074242137349
$ws.Range("F23").Value = 31.49351879877965
$ws.Range("G23").Value = 41.39349869852795
$ws.Range("H23").Value = 17.07939263160949
$ws.Range("K23").Value = 11.66504865222546
$ws.Range("L23").Value = 10.97010601064638
$ws.Range("N23").Value = 20.27510951975506
$ws.Range("B24").Value = 15.67171593495422
$ws.Range("C24").Value = 8.961351927349389
$ws.Range("D24").Value = 6.245946459655029
$ws.Range("F24").Value = 31.34682017208896
$ws.Range("G24").Value = 41.12964065485961
$ws.Range("H24").Value = 17.10683553176984
$ws.Range("K24").Value = 11.36135489319101
$ws.Range("L24").Value = 10.90124823263696
$ws.Range("N24").Value = 20.36870693070652
$ws.Range("B25").Value = 15.192620054445
$ws.Range("C25").Value = 8.940945386051881
$ws.Range("D25").Value = 6.22237252206478
$ws.Range("F25").Value = 31.21680483628908
$ws.Range("G25").Value = 40.88602109189537
$ws.Range("H25").Value = 17.14747102295115
$ws.Range("K25").Value = 11.0335724448885
$ws.Range("L25").Value = 10.8347996427775
$ws.Range("N25").Value = 20.47660683434466
